$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new questions text. The text literally contains the escaped
# sequence \u2019 (backslash, u, 2, 0, 1, 9) as plain characters rather than
# an actual Unicode right single quote character, matching the target file.
$newText = @'
questions = [
    {
        "title": "Sally\u2019s Salsa Shop has a balance sheet value of equities and liabilities of $10m. It is acquired by Tim\u2019s Taco Truck for $15m. How is the extra $5m recorded on the balance sheet?",
        "ques_type": 2,
        "options": [
            "As sunk costs",
            "As goodwill",
            "As non-current liabilities",
            "It is not recorded on the balance sheet"
        ],
        "score": "As goodwill"
    },
    {
        "title": "Which valuation method should be used for a company with uncertain cash flows in a volatile industry (e.g., a tech startup)?",
        "ques_type": 2,
        "options": [
            "A discounted cash flow: The uncertainty will be accounted for in the discount rate.",
            "Precedent transactions: This shows how previous investors have assessed the risk of the uncertain cashflows.",
            "Net income valuation: This shows you how much income the company has actually made, which is a good indication of the company\u2019s future earnings.",
            "Asset valuation: Looking at a company\u2019s assets is the only tangible way to value a company."
        ],
        "score": "Precedent transactions: This shows how previous investors have assessed the risk of the uncertain cashflows."
    },
    {
        "title": "A financial model workbook contains the following tabs. What kind of tabs are these?",
        "ques_type": 2,
        "options": [
            "Calculation tabs",
            "Input tabs",
            "Accounting tabs",
            "Model tabs"
        ],
        "score": "Calculation tabs"
    },
    {
        "title": "Why would you use this view of formulas?",
        "ques_type": 2,
        "options": [
            "To find inconsistencies or hard-coded inputs.",
            "To trace dependents and precedents.",
            "To set the cell formats.",
            "To indicate cells that still require changes."
        ],
        "score": "To find inconsistencies or hard-coded inputs."
    }
]
'@

# Remove trailing newline introduced by the here-string terminator.
$newText = $newText.TrimEnd("`r", "`n")

# Remove the old row 2 entirely (it held the shared string value) and clear
# any special formatting (bold font / border / alignment) from A1.
$ws.Rows.Item(2).Delete()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText
$ws.Rows.Item(1).AutoFit()
